# Refresh the crypto symbol list snapshot (Price / Volume(1h) columns)
# to match the GitHub Actions scrape taken on 2023-01-21 15:26 UTC.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store numeric-looking values
# as plain text (e.g. "303.46", "4.91%"). Force the Text number format
# on exactly the cells being refreshed (one call per contiguous block,
# since a multi-area union only applies NumberFormat to its first area)
# so Excel keeps writing them as text instead of auto-converting them
# to a number/percentage.
$ws.Range("D2:D9").NumberFormat = "@"
$ws.Range("D11:D27").NumberFormat = "@"
$ws.Range("D39:D51").NumberFormat = "@"
$ws.Range("E2:E27").NumberFormat = "@"
$ws.Range("E39:E51").NumberFormat = "@"

# Row 2: D2=303.46, E2=4.91%
$ws.Range("D2").Value = "303.46"
$ws.Range("E2").Value = "4.91%"

# Row 3: D3=34.92, E3=12.60%
$ws.Range("D3").Value = "34.92"
$ws.Range("E3").Value = "12.60%"

# Row 4: D4=5.162, E4=4.31%
$ws.Range("D4").Value = "5.162"
$ws.Range("E4").Value = "4.31%"

# Row 5: D5=0.07858, E5=6.68%
$ws.Range("D5").Value = "0.07858"
$ws.Range("E5").Value = "6.68%"

# Row 6: D6=2.329, E6=1.56%
$ws.Range("D6").Value = "2.329"
$ws.Range("E6").Value = "1.56%"

# Row 7: D7=8.055, E7=5.10%
$ws.Range("D7").Value = "8.055"
$ws.Range("E7").Value = "5.10%"

# Row 8: D8=3.976, E8=6.37%
$ws.Range("D8").Value = "3.976"
$ws.Range("E8").Value = "6.37%"

# Row 9: D9=0.9263, E9=1.05%
$ws.Range("D9").Value = "0.9263"
$ws.Range("E9").Value = "1.05%"

# Row 10: E10=10.37%
$ws.Range("E10").Value = "10.37%"

# Row 11: D11=0.1824, E11=7.07%
$ws.Range("D11").Value = "0.1824"
$ws.Range("E11").Value = "7.07%"

# Row 12: D12=0.08531, E12=3.16%
$ws.Range("D12").Value = "0.08531"
$ws.Range("E12").Value = "3.16%"

# Row 13: D13=0.03407, E13=9.39%
$ws.Range("D13").Value = "0.03407"
$ws.Range("E13").Value = "9.39%"

# Row 14: D14=0.09910, E14=-0.78%
$ws.Range("D14").Value = "0.09910"
$ws.Range("E14").Value = "-0.78%"

# Row 15: D15=0.001477, E15=-1.48%
$ws.Range("D15").Value = "0.001477"
$ws.Range("E15").Value = "-1.48%"

# Row 16: D16=0.005815, E16=1.13%
$ws.Range("D16").Value = "0.005815"
$ws.Range("E16").Value = "1.13%"

# Row 17: D17=3.471, E17=-0.02%
$ws.Range("D17").Value = "3.471"
$ws.Range("E17").Value = "-0.02%"

# Row 18: D18=2.103, E18=-1.01%
$ws.Range("D18").Value = "2.103"
$ws.Range("E18").Value = "-1.01%"

# Row 19: D19=0.3429, E19=3.01%
$ws.Range("D19").Value = "0.3429"
$ws.Range("E19").Value = "3.01%"

# Row 20: D20=0.1325, E20=1.96%
$ws.Range("D20").Value = "0.1325"
$ws.Range("E20").Value = "1.96%"

# Row 21: D21=4.532, E21=8.66%
$ws.Range("D21").Value = "4.532"
$ws.Range("E21").Value = "8.66%"

# Row 22: D22=0.2219, E22=4.58%
$ws.Range("D22").Value = "0.2219"
$ws.Range("E22").Value = "4.58%"

# Row 23: D23=0.04633, E23=2.87%
$ws.Range("D23").Value = "0.04633"
$ws.Range("E23").Value = "2.87%"

# Row 24: D24=0.001215, E24=0.09%
$ws.Range("D24").Value = "0.001215"
$ws.Range("E24").Value = "0.09%"

# Row 25: D25=0.004459, E25=6.22%
$ws.Range("D25").Value = "0.004459"
$ws.Range("E25").Value = "6.22%"

# Row 26: D26=0.0001296, E26=-0.11%
$ws.Range("D26").Value = "0.0001296"
$ws.Range("E26").Value = "-0.11%"

# Row 27: D27=0.0003391, E27=-0.01%
$ws.Range("D27").Value = "0.0003391"
$ws.Range("E27").Value = "-0.01%"

# Row 39: D39=0.01756, E39=11.76%
$ws.Range("D39").Value = "0.01756"
$ws.Range("E39").Value = "11.76%"

# Row 40: D40=0.04737, E40=5.26%
$ws.Range("D40").Value = "0.04737"
$ws.Range("E40").Value = "5.26%"

# Row 41: D41=0.007772, E41=5.38%
$ws.Range("D41").Value = "0.007772"
$ws.Range("E41").Value = "5.38%"

# Row 42: D42=0.1415, E42=5.91%
$ws.Range("D42").Value = "0.1415"
$ws.Range("E42").Value = "5.91%"

# Row 43: D43=0.008803, E43=-10.58%
$ws.Range("D43").Value = "0.008803"
$ws.Range("E43").Value = "-10.58%"

# Row 44: D44=0.002284, E44=3.04%
$ws.Range("D44").Value = "0.002284"
$ws.Range("E44").Value = "3.04%"

# Row 45: D45=0.009169, E45=7.57%
$ws.Range("D45").Value = "0.009169"
$ws.Range("E45").Value = "7.57%"

# Row 46: D46=0.00006051, E46=-0.85%
$ws.Range("D46").Value = "0.00006051"
$ws.Range("E46").Value = "-0.85%"

# Row 47: D47=0.00000000748, E47=-0.12%
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").Value = "-0.12%"

# Row 48: D48=3.886, E48=49.11%
$ws.Range("D48").Value = "3.886"
$ws.Range("E48").Value = "49.11%"

# Row 49: D49=0.002684, E49=34.33%
$ws.Range("D49").Value = "0.002684"
$ws.Range("E49").Value = "34.33%"

# Row 50: D50=0.00002095, E50=-0.12%
$ws.Range("D50").Value = "0.00002095"
$ws.Range("E50").Value = "-0.12%"

# Row 51: D51=0.0001996, E51=-0.12%
$ws.Range("D51").Value = "0.0001996"
$ws.Range("E51").Value = "-0.12%"
